$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.444.26"
$ws.Range("E2").Value = "  -5.61%  "
$ws.Range("D3").Value = "2.616.40"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'301.55"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").Value = "'95.59"
$ws.Range("E6").Value = "  -3.74%  "
$ws.Range("D7").Value = "'0.580"
$ws.Range("E7").Value = "  -3.58%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.553"
$ws.Range("E9").Value = "  -4.17%  "
$ws.Range("D10").Value = "'36.88"
$ws.Range("E10").Value = "  -5.79%  "
$ws.Range("D11").Value = "'0.0810"
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("D12").Value = "'7.76"
$ws.Range("E12").Value = "  -4.95%  "
$ws.Range("D13").Value = "3.031.78"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "2.653.91"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "'0.885"
$ws.Range("E16").Value = "  -3.17%  "
$ws.Range("D17").Value = "'14.41"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").Value = "43.669.00"
$ws.Range("E18").Value = "  -5.80%  "
$ws.Range("D19").Value = "'6.64"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").Value = "0.0₃0970"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("D21").Value = "'12.39"
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("D22").Value = "'73.48"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("D23").Value = "'266.98"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.93"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'2.21"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("D26").Value = "'29.50"
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "'10.19"
$ws.Range("E28").Value = "  -3.03%  "
$ws.Range("E29").Value = "  -3.28%  "
$ws.Range("D30").Value = "'37.47"
$ws.Range("E30").Value = "  -4.41%  "
$ws.Range("D31").Value = "'6.03"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").Value = "'3.60"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "'2.23"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("B34").Value = "EnergySwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D34").Value = "'28.06"
$ws.Range("E34").Value = "  +20.41%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'152.56"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.79"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.0808"
$ws.Range("E37").Value = "  -3.36%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.118"
$ws.Range("E38").Value = "  -3.05%  "
$ws.Range("D39").Value = "'0.120"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Value = "'15.83"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "'3.48"
$ws.Range("E41").Value = "  -3.28%  "
$ws.Range("D42").Value = "'0.0314"
$ws.Range("E42").Value = "  -3.96%  "
$ws.Range("D43").Value = "'3.82"
$ws.Range("E43").Value = "  -6.09%  "
$ws.Range("D44").Value = "2.082.05"
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").Value = "'88.64"
$ws.Range("E46").Value = "  -4.64%  "
$ws.Range("D47").Value = "'9.10"
$ws.Range("E47").Value = "  -5.94%  "
$ws.Range("D48").Value = "2.892.07"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'1.59"
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'106.29"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").Value = "'0.190"
$ws.Range("E51").Value = "  -4.34%  "
